# Chronomancer.xlsx — "Add 3 new relics"
#
# The only substantive content change in the target revision lives on the
# "Relic" worksheet: several relic effects were reworded, the "Goggles"
# relic was re-rarified from Rare to Uncommon, three brand new relics were
# added ("Grimoire", "Winder", "Old Wine") and the table shrank from 15 to
# 14 relic rows (16 sheet rows -> 15 sheet rows, including the header).
# Every other worksheet in the workbook only shows shared-string index
# renumbering caused by this edit, with no actual value changes, so
# nothing else needs to be touched (Excel recomputes the shared string
# table automatically when the file is saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Relic")

# The table used to have 15 relic rows (2-16); it now has 14 (2-15).
# Remove the old trailing row so everything below shifts up and the
# dimension / trailing formatting rows line up with the target file.
$ws.Rows.Item(16).Delete()

# Row 2 - Broken Watch (Starting) - effect text unchanged, now marked Done.
$ws.Range("A2").Value = "Broken Watch"
$ws.Range("B2").Value = "Starting"
$ws.Range("C2").Value = "The first time you consume Jade each combat, gain E."
$ws.Range("D2").Value = "Yes"
$ws.Range("E2").Value = $true

# Row 3 - Common relic #1 - reworded effect.
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = "Common"
$ws.Range("C3").Value = "Whenever you trigger a Recall effect, deal 2 damage to a random enemy."
$ws.Range("D3").Value = "Yes"
$ws.Range("E3").ClearContents()

# Row 4 - Common relic #2 - reworded effect; row height shrinks back to single line.
$ws.Range("A4").ClearContents()
$ws.Range("B4").Value = "Common"
$ws.Range("C4").Value = "Whenever you play 4 cards that cost 0, draw 1 card."
$ws.Range("D4").Value = "No"
$ws.Range("E4").ClearContents()
$ws.Rows.Item(4).RowHeight = 17

# Row 5 - Common relic #3 - reworded effect.
$ws.Range("A5").ClearContents()
$ws.Range("B5").Value = "Common"
$ws.Range("C5").Value = "At the start of the battle, you can choose to place up to 4 cards into your discard pile."
$ws.Range("D5").Value = "No"
$ws.Range("E5").ClearContents()

# Row 6 - Uncommon relic #1 - "increased by 4" -> "increased by 3".
$ws.Range("A6").ClearContents()
$ws.Range("B6").Value = "Uncommon"
$ws.Range("C6").Value = "The first cost X card you play each combat has its effect increased by 3."
$ws.Range("D6").Value = "Yes"
$ws.Range("E6").ClearContents()

# Row 7 - Uncommon relic #2 - reworded effect.
$ws.Range("A7").ClearContents()
$ws.Range("B7").Value = "Uncommon"
$ws.Range("C7").Value = "Whenever you consume Jade, gain 5 Block."
$ws.Range("D7").Value = "Yes"
$ws.Range("E7").ClearContents()

# Row 8 - Goggles moves here from the old Rare slot, re-rarified to
# Uncommon, no longer Exclusive, and marked Done; row height shrinks back
# to single line.
$ws.Range("A8").Value = "Goggles"
$ws.Range("B8").Value = "Uncommon"
$ws.Range("C8").Value = "The first card you play each combat costs 0."
$ws.Range("D8").Value = "No"
$ws.Range("E8").Value = $true
$ws.Rows.Item(8).RowHeight = 17

# Row 9 - NEW relic: Grimoire (Rare).
$ws.Range("A9").Value = "Grimoire"
$ws.Range("B9").Value = "Rare"
$ws.Range("C9").Value = "Whenever you spend E on a cost X card, a random card in your hand costs 0 this turn."
$ws.Range("D9").Value = "Yes"
$ws.Range("E9").Value = $true

# Row 10 - Rare relic - effect text unchanged.
$ws.Range("A10").ClearContents()
$ws.Range("B10").Value = "Rare"
$ws.Range("C10").Value = "At the end of your turn, if you have 3 or more cards in your hand, Recall: Draw 2 cards."
$ws.Range("D10").Value = "Yes"
$ws.Range("E10").ClearContents()

# Row 11 - Rare relic - effect text unchanged.
$ws.Range("A11").ClearContents()
$ws.Range("B11").Value = "Rare"
$ws.Range("C11").Value = "At the start of your third turn, remove all your debuffs."
$ws.Range("D11").Value = "No"
$ws.Range("E11").ClearContents()

# Row 12 - Ancient Watch (Boss) - effect text unchanged, now marked Done.
$ws.Range("A12").Value = "Ancient Watch"
$ws.Range("B12").Value = "Boss"
$ws.Range("C12").Value = "The first 3 times you consume Jade each combat, gain E."
$ws.Range("D12").Value = "Yes"
$ws.Range("E12").Value = $true

# Row 13 - Boss relic - "Remove 1 Jade" -> "Remove a Jade".
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = "Boss"
$ws.Range("C13").Value = "Remove a Jade at the end of your turn."
$ws.Range("D13").Value = "Yes"
$ws.Range("E13").ClearContents()

# Row 14 - NEW relic: Winder (Boss).
$ws.Range("A14").Value = "Winder"
$ws.Range("B14").Value = "Boss"
$ws.Range("C14").Value = "The first time you trigger a Recall effect each turn, gain E."
$ws.Range("D14").Value = "Yes"
$ws.Range("E14").Value = $true

# Row 15 - NEW relic: Old Wine (Shop) - effect text unchanged from the old
# unnamed Shop relic, now named and marked Done.
$ws.Range("A15").Value = "Old Wine"
$ws.Range("B15").Value = "Shop"
$ws.Range("C15").Value = "At the start of your turn, gain 2 Strength and 1 Jade."
$ws.Range("D15").Value = "Yes"
$ws.Range("E15").Value = $true

# Match the author's final cursor position recorded in the sheet view.
$ws.Activate()
$ws.Range("E16").Select()
